# fix: errors in spelling
# The cell B2 on the "Geofence Visits" sheet contained the German text
# "Geofence besuche" which should read "Besuchte Orte".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Geofence Visits")

$ws.Range("B2").Value = "Besuchte Orte"

# Reflect the selection change recorded alongside the edit (user selected B2).
$ws.Range("B2").Select()
